# Generate Report for Handoff
# Adds a new report row (for 415b514f-c819-465e-b750-d421ae11c289.md) to the
# Overview sheet and to each locale sheet (zh-cn, de-de), mirroring the
# existing 071d1400-... row that was generated for the first handoff.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTest/oltest/blob/6feb61ca7a8b9c94737b3eab6a938bf85850b170/e2e/"
$newFile = "415b514f-c819-465e-b750-d421ae11c289.md"
$newFileUrl = $baseUrl + $newFile

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-40-13 06:40:02"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $newFileUrl, "", "", $newFile)

# ---------------------------------------------------------------------------
# zh-cn sheet: Source File Name | File Extension | Status | Latest Handoff
# File | Latest Handoff Datetime | ... | Latest Handback DateTime |
# Handoff Reason
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhTarget = "415b514f-c819-465e-b750-d421ae11c289.5e2c28fec218bce00dcb0c8195fd0f47ab8ef283.zh-cn.xlf"
$zhTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f6078494afc9cc790b0120a8ee67269ae5b15e79/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $zhTarget

$wsZh.Range("A3").Value = $newFile
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = $zhTarget
$wsZh.Range("E3").Value = "2016-03-13 06:39:58"
$wsZh.Range("H3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I3").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $newFileUrl, "", "", $newFile)
$wsZh.Hyperlinks.Add($wsZh.Range("B3"), $newFileUrl, "", "", ".md")
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $zhTargetUrl, "", "", $zhTarget)

# ---------------------------------------------------------------------------
# de-de sheet: same shape as zh-cn
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deTarget = "415b514f-c819-465e-b750-d421ae11c289.5e2c28fec218bce00dcb0c8195fd0f47ab8ef283.de-de.xlf"
$deTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6b0d4fd2a6a03fa1dea586f917bc8ac91291ecf5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $deTarget

$wsDe.Range("A3").Value = $newFile
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = $deTarget
$wsDe.Range("E3").Value = "2016-03-13 06:40:02"
$wsDe.Range("H3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I3").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $newFileUrl, "", "", $newFile)
$wsDe.Hyperlinks.Add($wsDe.Range("B3"), $newFileUrl, "", "", ".md")
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $deTargetUrl, "", "", $deTarget)
